$wb = $excel.ActiveWorkbook

# --- Alternative sheet (sheet1) ---
$wsAlt = $wb.Worksheets.Item("Alternative")
$wsAlt.Range("N2").Value = -10.57
$wsAlt.Range("N3").Value = -4.81
$wsAlt.Range("N4").Value = -9.289999999999999
$wsAlt.Range("N5").Value = -82.09999999999999
$wsAlt.Range("N6").Value = -8.539999999999999
$wsAlt.Range("N7").Value = -7.38

# --- Bond sheet (sheet2) ---
$wsBond = $wb.Worksheets.Item("Bond")
$wsBond.Range("N2").Value = -6.37
$wsBond.Range("N3").Value = 3.69
$wsBond.Range("N4").Value = -26.16
$wsBond.Range("N5").Value = -16.33

# --- Equity sheet (sheet3) ---
$wsEq = $wb.Worksheets.Item("Equity")
$wsEq.Range("N2").Value = -17.55
$wsEq.Range("N3").Value = -15.7
$wsEq.Range("N4").Value = -18.5
$wsEq.Range("N5").Value = -8.890000000000001
$wsEq.Range("N6").Value = -37.84
$wsEq.Range("N7").Value = -17.07
$wsEq.Range("N8").Value = -19.12
$wsEq.Range("N9").Value = -23.39
$wsEq.Range("N10").Value = -35.39
$wsEq.Range("N11").Value = -1.31
